# Commit: "Add data for 2024-11-15"
# The sheet tracks one day of readings per column (B..BO already hold
# 2024/09/09 .. 2024/11/14). This appends a new trailing column BP with the
# 2024/11/15 readings: a text date header in row 1 and one numeric value per
# machine in rows 2-53. Some values are highlighted (yellow / light-blue
# fill) exactly like the existing columns, so we copy that cell formatting
# across from an already-highlighted cell instead of re-describing it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cells already carrying each of the three cell styles used on
# this sheet: 1 = plain, 2 = yellow fill, 3 = light-blue fill.
$styleSource = @{ 1 = "A2"; 2 = "BM2"; 3 = "BC2" }

# --- Row 1: new date header, kept as literal text like the other headers ---
$ws.Range("BP1").NumberFormat = "@"
$ws.Range("BP1").Value = "2024/11/15"
$ws.Range("BO1").Copy()
$ws.Range("BP1").PasteSpecial(-4122)

# --- Rows 2-53: the 2024/11/15 reading for every machine ---
$newData = @(
    [PSCustomObject]@{ Row = 2; Value = 163.1; Style = 1 }
    [PSCustomObject]@{ Row = 3; Value = 148.7; Style = 1 }
    [PSCustomObject]@{ Row = 4; Value = 191.4; Style = 1 }
    [PSCustomObject]@{ Row = 5; Value = 146.5; Style = 1 }
    [PSCustomObject]@{ Row = 6; Value = 145.3; Style = 1 }
    [PSCustomObject]@{ Row = 7; Value = 142.5; Style = 1 }
    [PSCustomObject]@{ Row = 8; Value = 165.8; Style = 1 }
    [PSCustomObject]@{ Row = 9; Value = 137.9; Style = 3 }
    [PSCustomObject]@{ Row = 10; Value = 122.1; Style = 2 }
    [PSCustomObject]@{ Row = 11; Value = 143.4; Style = 1 }
    [PSCustomObject]@{ Row = 12; Value = 142; Style = 1 }
    [PSCustomObject]@{ Row = 13; Value = 233.5; Style = 1 }
    [PSCustomObject]@{ Row = 14; Value = 149.7; Style = 1 }
    [PSCustomObject]@{ Row = 15; Value = 118.6; Style = 2 }
    [PSCustomObject]@{ Row = 16; Value = 144.1; Style = 1 }
    [PSCustomObject]@{ Row = 17; Value = 131; Style = 3 }
    [PSCustomObject]@{ Row = 18; Value = 157.2; Style = 1 }
    [PSCustomObject]@{ Row = 19; Value = 185.5; Style = 1 }
    [PSCustomObject]@{ Row = 20; Value = 164.4; Style = 1 }
    [PSCustomObject]@{ Row = 21; Value = 166.9; Style = 1 }
    [PSCustomObject]@{ Row = 22; Value = 117.8; Style = 2 }
    [PSCustomObject]@{ Row = 23; Value = 130.8; Style = 3 }
    [PSCustomObject]@{ Row = 24; Value = 163.5; Style = 1 }
    [PSCustomObject]@{ Row = 25; Value = 166.5; Style = 1 }
    [PSCustomObject]@{ Row = 26; Value = 141.6; Style = 1 }
    [PSCustomObject]@{ Row = 27; Value = 155.9; Style = 1 }
    [PSCustomObject]@{ Row = 28; Value = 153.1; Style = 1 }
    [PSCustomObject]@{ Row = 29; Value = 138.9; Style = 3 }
    [PSCustomObject]@{ Row = 30; Value = 145.1; Style = 1 }
    [PSCustomObject]@{ Row = 31; Value = 155.8; Style = 1 }
    [PSCustomObject]@{ Row = 32; Value = 140.6; Style = 1 }
    [PSCustomObject]@{ Row = 33; Value = 120.6; Style = 2 }
    [PSCustomObject]@{ Row = 34; Value = 192.1; Style = 1 }
    [PSCustomObject]@{ Row = 35; Value = 153.4; Style = 1 }
    [PSCustomObject]@{ Row = 36; Value = 149.1; Style = 1 }
    [PSCustomObject]@{ Row = 37; Value = 209; Style = 1 }
    [PSCustomObject]@{ Row = 38; Value = 131.8; Style = 3 }
    [PSCustomObject]@{ Row = 39; Value = 134.3; Style = 3 }
    [PSCustomObject]@{ Row = 40; Value = 121.9; Style = 2 }
    [PSCustomObject]@{ Row = 41; Value = 182.8; Style = 1 }
    [PSCustomObject]@{ Row = 42; Value = 202.1; Style = 1 }
    [PSCustomObject]@{ Row = 43; Value = 164.8; Style = 1 }
    [PSCustomObject]@{ Row = 44; Value = 146.5; Style = 1 }
    [PSCustomObject]@{ Row = 45; Value = 114.5; Style = 2 }
    [PSCustomObject]@{ Row = 46; Value = 159.6; Style = 1 }
    [PSCustomObject]@{ Row = 47; Value = 152.5; Style = 1 }
    [PSCustomObject]@{ Row = 48; Value = 147; Style = 1 }
    [PSCustomObject]@{ Row = 49; Value = 164.8; Style = 1 }
    [PSCustomObject]@{ Row = 50; Value = 152.8; Style = 1 }
    [PSCustomObject]@{ Row = 51; Value = 144.6; Style = 1 }
    [PSCustomObject]@{ Row = 52; Value = 144.7; Style = 1 }
    [PSCustomObject]@{ Row = 53; Value = 125.8; Style = 3 }
)

foreach ($item in $newData) {
    $cell = "BP$($item.Row)"
    $ws.Range($cell).Value = $item.Value
    $ws.Range($styleSource[$item.Style]).Copy()
    $ws.Range($cell).PasteSpecial(-4122)
}

$excel.CutCopyMode = $false

# New column BP needs the same "12" width as every other data column.
$ws.Columns.Item(68).ColumnWidth = 12 - 5/6
